$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write one full data row (columns B..AC) into a given sheet row,
# leaving column A (the running "id" index) untouched.
# ---------------------------------------------------------------------------
function Set-RowData {
    param($sheet, [int]$row, [hashtable]$vals)

    foreach ($col in $vals.Keys) {
        $sheet.Range("$col$row").Value = $vals[$col]
    }
}

# swap 130 <-> 131
$row130 = @{ B=7453204; E=45242.79166666666; F='Cerro Porteno'; G='Tacuary'; H=1; I=1; J='D'; K=1.285; L=5; M=8; N=1.285; O=4.75; P=8; Q=-1.5; R=1.9; S=1.9; T=3; U=1.9; V=1.9; W=-1; X=3.75; Y=-1; Z=-1; AA=0.8999999999999999; AB=-1; AC=0.8999999999999999 }
$row131 = @{ B=7454842; E=45242.79166666666; F='Sportivo Luqueno'; G='Libertad Asuncion'; H=0; I=1; J='A'; K=4; L=3.6; M=1.727; N=3.5; O=3.3; P=1.95; Q=0.5; R=1.8; S=2; T=2.5; U=1.975; V=1.825; W=-1; X=-1; Y=0.95; Z=-1; AA=1; AB=-1; AC=0.825 }
Set-RowData $ws 130 $row130
Set-RowData $ws 131 $row131

# swap 134 <-> 135
$row134 = @{ B=7493428; E=45253.8125; F='Guairena FC'; G='Resistencia FC'; H=4; I=1; J='H'; K=1.727; L=3.6; M=4.2; N=1.45; O=4.2; P=6; Q=-1; R=1.775; S=2.025; T=2.75; U=1.825; V=1.975; W=0.45; X=-1; Y=-1; Z=0.7749999999999999; AA=-1; AB=0.825; AC=-1 }
$row135 = @{ B=7493427; E=45253.8125; F='Tacuary'; G='Sportivo Luqueno'; H=1; I=1; J='D'; K=3.4; L=3.3; M=2; N=3.2; O=3.25; P=2.1; Q=0.25; R=2.025; S=1.775; T=2.5; U=1.975; V=1.825; W=-1; X=2.25; Y=-1; Z=0.5125; AA=-0.5; AB=-1; AC=0.825 }
Set-RowData $ws 134 $row134
Set-RowData $ws 135 $row135

# swap 137 <-> 138
$row137 = @{ B=7499442; E=45255.8125; F='Guarani Asuncion'; G='Sportivo Trinidense'; H=0; I=2; J='A'; K=2; L=3.4; M=3.3; N=1.909; O=3.4; P=3.6; Q=-0.5; R=1.975; S=1.825; T=2.5; U=1.9; V=1.9; W=-1; X=-1; Y=2.6; Z=-1; AA=0.825; AB=-1; AC=0.8999999999999999 }
$row138 = @{ B=7499440; E=45255.8125; F='Olimpia Asuncion'; G='Libertad Asuncion'; H=1; I=3; J='A'; K=2.7; L=3.4; M=2.3; N=2.8; O=3.25; P=2.3; Q=0.25; R=1.75; S=2.05; T=2.25; U=1.85; V=1.95; W=-1; X=-1; Y=1.3; Z=-1; AA=1.05; AB=0.8500000000000001; AC=-1 }
Set-RowData $ws 137 $row137
Set-RowData $ws 138 $row138

# swap 140 <-> 141
$row140 = @{ B=7493310; E=45259.8125; F='Libertad Asuncion'; G='Tacuary'; H=1; I=2; J='A'; K=1.363; L=5; M=7; N=1.571; O=4.2; P=4.75; Q=-0.75; R=1.8; S=2; T=2.75; U=1.8; V=2; W=-1; X=-1; Y=3.75; Z=-1; AA=1; AB=0.4; AC=-0.5 }
$row141 = @{ B=7493431; E=45259.8125; F='Sportivo Trinidense'; G='Guairena FC'; H=7; I=2; J='H'; K=2.05; L=3.3; M=3.3; N=2.6; O=3.1; P=2.6; Q=0; R=1.925; S=1.875; T=2.5; U=2; V=1.8; W=1.6; X=-1; Y=-1; Z=0.925; AA=-1; AB=1; AC=-1 }
Set-RowData $ws 140 $row140
Set-RowData $ws 141 $row141

# swap 143 <-> 145
$row143 = @{ B=7493312; E=45261.8125; F='Cerro Porteno'; G='Guarani Asuncion'; H=4; I=0; J='H'; K=1.7; L=3.6; M=4.333; N=1.727; O=3.75; P=4.2; Q=-0.5; R=1.8; S=2; T=2.75; U=1.875; V=1.925; W=0.7270000000000001; X=-1; Y=-1; Z=0.8; AA=-1; AB=0.875; AC=-1 }
$row145 = @{ B=7493433; E=45261.8125; F='Sportivo Luqueno'; G='Nacional Asuncion'; H=1; I=1; J='D'; K=2.75; L=3.2; M=2.4; N=2.75; O=3.1; P=2.45; Q=0.25; R=1.75; S=2.05; T=2.25; U=2; V=1.8; W=-1; X=2.1; Y=-1; Z=0.375; AA=-0.5; AB=-0.5; AC=0.4 }
Set-RowData $ws 143 $row143
Set-RowData $ws 145 $row145

# ---------------------------------------------------------------------------
# Row 206: this fixture's final result / closing-odds data has come in.
# It now carries the match that used to be parked at row 209 (2 de Mayo
# vs Guarani Asuncion), with final closing odds filled in.
# ---------------------------------------------------------------------------
$row206 = @{ B=7609672; E=45382.875; F='2 de Mayo'; G='Guarani Asuncion'; K=2.2; L=3.1; M=3.1; N=2.55; O=3; P=2.7; Q=0; R=1.85; S=1.95; T=2.25; U=1.975; V=1.825; W=0; X=0; Y=0; Z=0; AA=0 }
Set-RowData $ws 206 $row206

# ---------------------------------------------------------------------------
# Rows 207-209 no longer exist in the refreshed extract - drop them.
# ---------------------------------------------------------------------------
$ws.Rows("207:209").Delete()
